# Weekly fruit/vegetable price update: a new weekly record is inserted
# into the "Perejil" (parsley) price log at row 73, pushing all
# subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 73 (shifts existing rows 73:106 down to 74:107,
# and copies formatting - including the date style on column D - from the row above).
$ws.Rows.Item(73).Insert()

# Populate the new row 73 with the new weekly price observation.
$ws.Range("A73").Value = 7
$ws.Range("B73").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C73").Value = "Ñuble"
$ws.Range("D73").Value = 45134
$ws.Range("E73").Value = 16
$ws.Range("F73").Value = 100112044
$ws.Range("G73").Value = "Perejil"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 100
$ws.Range("K73").Value = 1500
$ws.Range("L73").Value = 1500
$ws.Range("M73").Value = 1500
$ws.Range("N73").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O73").Value = "Región de Ñuble"
$ws.Range("P73").Value = 1500
$ws.Range("Q73").Value = 1
$ws.Range("R73").Value = "Hortaliza"
